$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 87; existing rows 87-120 shift down to 88-121.
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new weekly record.
$ws.Cells.Item(87, 1).Value  = 8
$ws.Cells.Item(87, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(87, 3).Value  = "Coquimbo"
$ws.Cells.Item(87, 4).Value  = 44463
$ws.Cells.Item(87, 5).Value  = 4
$ws.Cells.Item(87, 6).Value  = 100112031
$ws.Cells.Item(87, 7).Value  = "Poroto verde"
$ws.Cells.Item(87, 8).Value  = "Magnum"
$ws.Cells.Item(87, 9).Value  = "Primera"
$ws.Cells.Item(87, 10).Value = 500
$ws.Cells.Item(87, 11).Value = 35000
$ws.Cells.Item(87, 12).Value = 36000
$ws.Cells.Item(87, 13).Value = 35500
$ws.Cells.Item(87, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(87, 15).Value = "Perú"
$ws.Cells.Item(87, 16).Value = 1420
$ws.Cells.Item(87, 17).Value = 25
$ws.Cells.Item(87, 18).Value = "Hortaliza"
